# Adiciona os dados da sala (turma) na caixa de texto da disciplina,
# no primeiro slide (slide de rosto / capa).
#
# A caixa de texto "CaixaDeTexto 16" contem hoje apenas um paragrafo com
# "DISCIPLINA:   COMPLIANCE & QUALITY ASSURANCE". O pedido ("adicionando
# dados da sala") acrescenta um novo paragrafo abaixo com o texto "3SI"
# (a turma/sala), em negrito, igual ao restante do nome da disciplina.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("CaixaDeTexto 16")

$tr = $sh.TextFrame.TextRange

# Insere uma quebra de paragrafo seguida do texto "3SI" logo apos o
# conteudo existente. A caixa tem spAutoFit habilitado, entao a altura
# da forma e recalculada automaticamente para acomodar a nova linha de
# texto (crescendo de 369332 para 646331 EMU).
[void]$tr.InsertAfter([char]13 + "3SI")
